# Bug fix: the "get_impact_data" test case row was a leftover duplicate;
# rename row 2 to "search_impact" (keeping its "Pass" status) and remove
# the now-redundant row 3 that used to hold "search_impact" / "Pass".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2, column A: get_impact_data -> search_impact (status in B2 stays "Pass")
$ws.Range("A2").Value = "search_impact"

# Remove the old row 3 (search_impact / Pass), shifting nothing else up
# since it was the last row; this also shrinks the used range to A1:B2.
$ws.Rows("3:3").Delete()
